# Append new tracker rows (72-83) to Sheet1, mirroring the existing
# alternating G1/G2 pattern for the next six calendar days.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date format already used by column C (e.g. C71) - reapply to the new
# date cells so they keep the YYYY-MM-DD format.
$dateFormat = $ws.Range("C71").NumberFormat

$rows = @(
    @{ Row = 72; A = "G1"; B = "Test1";              C = 45896; D = 0.7129733411147666; E = 0; F = -0.01 },
    @{ Row = 73; A = "G1"; B = "Test1";              C = 45897; D = 0.7059141991235313; E = 0; F = -0.01 },
    @{ Row = 74; A = "G1"; B = "Test1";              C = 45898; D = 0.6989249496272587; E = 0; F = -0.01 },
    @{ Row = 75; A = "G1"; B = "Test1";              C = 45899; D = 0.6920049006210482; E = 0; F = -0.01 },
    @{ Row = 76; A = "G1"; B = "Test1";              C = 45900; D = 0.6851533669515329; E = 0; F = -0.01 },
    @{ Row = 77; A = "G1"; B = "Test1";              C = 45901; D = 0.6783696702490425; E = 0; F = -0.01 },
    @{ Row = 78; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45896; D = 0.7129733411147666; E = 0; F = -0.01 },
    @{ Row = 79; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45897; D = 0.7059141991235313; E = 0; F = -0.01 },
    @{ Row = 80; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45898; D = 0.6989249496272587; E = 0; F = -0.01 },
    @{ Row = 81; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45899; D = 0.6920049006210482; E = 0; F = -0.01 },
    @{ Row = 82; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45900; D = 0.6851533669515329; E = 0; F = -0.01 },
    @{ Row = 83; A = "G2"; B = "sedrftgyhuioygtfrd"; C = 45901; D = 0.6783696702490425; E = 0; F = -0.01 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 3).NumberFormat = $dateFormat
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
